# Remove the stray duplicate "Resumo dos PADs" slide (5th slide in the deck).
# It duplicated the Right Arrow / Rounded Rectangle / Table layout used by the
# neighboring slides but pointed at an older snapshot image (image6.png),
# while the following slides already carry the corrected image (image7.png).
$p = $ppt.ActivePresentation
$p.Slides.Item(5).Delete()
